$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4066.3333
$ws.Range("I62").Value = 4066.3333
$ws.Range("K62").Value = 4066.3333
$ws.Range("M62").Value = -3442.3333

$ws.Range("H65").Value = 4066.3333
$ws.Range("I65").Value = 4066.3333
$ws.Range("K65").Value = 20331.6665
$ws.Range("M65").Value = -17211.6665

$ws.Range("H70").Value = 3940.2
$ws.Range("I70").Value = 7550.5
$ws.Range("K70").Value = 22651.5
$ws.Range("M70").Value = -22381.5

$ws.Range("H73").Value = 3940.2
$ws.Range("I73").Value = 7550.5
$ws.Range("K73").Value = 22651.5
$ws.Range("M73").Value = -21715.5

$ws.Range("H88").Value = 2486.5557
$ws.Range("I88").Value = 2895
$ws.Range("J88").Value = 2369.8572
$ws.Range("K88").Value = 2895
$ws.Range("L88").Value = 2369.8572
$ws.Range("M88").Value = -2489
$ws.Range("N88").Value = -3181.8572

$ws.Range("H91").Value = 2486.5557
$ws.Range("I91").Value = 2895
$ws.Range("J91").Value = 2369.8572
$ws.Range("K91").Value = 2895
$ws.Range("L91").Value = 2369.8572
$ws.Range("M91").Value = -1491
$ws.Range("N91").Value = -5177.8572

$ws.Range("H132").Value = 4316.48
$ws.Range("I132").Value = 4316.48
$ws.Range("K132").Value = 12949.44
$ws.Range("M132").Value = -10419.44

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 978.7368
$ws.Range("I2").Value = 1018.05884
$ws.Range("K2").Value = 1018.05884
$ws.Range("M2").Value = -905.05884

$ws.Range("H61").Value = 1613.8462
$ws.Range("I61").Value = 1613.8462
$ws.Range("K61").Value = 1613.8462
$ws.Range("M61").Value = -1401.8462

$ws.Range("H88").Value = 2061.6
$ws.Range("I88").Value = 1446.5
$ws.Range("J88").Value = 2215.375
$ws.Range("K88").Value = 1446.5
$ws.Range("L88").Value = 2215.375
$ws.Range("M88").Value = -1040.5
$ws.Range("N88").Value = -3027.375

$ws.Range("H91").Value = 2061.6
$ws.Range("I91").Value = 1446.5
$ws.Range("J91").Value = 2215.375
$ws.Range("K91").Value = 1446.5
$ws.Range("L91").Value = 2215.375
$ws.Range("M91").Value = -42.5
$ws.Range("N91").Value = -5023.375

$ws.Range("H97").Value = 2709.5
$ws.Range("I97").Value = 1877.5
$ws.Range("K97").Value = 1877.5
$ws.Range("M97").Value = -1381.5

$ws.Range("H116").Value = 978.7368
$ws.Range("I116").Value = 1018.05884
$ws.Range("K116").Value = 1018.05884
$ws.Range("M116").Value = 1275.94116

$ws.Range("H136").Value = 1613.8462
$ws.Range("I136").Value = 1613.8462
$ws.Range("K136").Value = 4841.5386
$ws.Range("M136").Value = -2291.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 978.7368
$ws.Range("I3").Value = 1018.05884
$ws.Range("K3").Value = 1018.05884
$ws.Range("M3").Value = -904.05884

$ws.Range("H20").Value = 3321.5715
$ws.Range("I20").Value = 2564.25
$ws.Range("K20").Value = 2564.25
$ws.Range("M20").Value = -2317.25

$ws.Range("H86").Value = 2261.926
$ws.Range("I86").Value = 2358.7727
$ws.Range("K86").Value = 2358.7727
$ws.Range("M86").Value = -1235.7727

$ws.Range("H89").Value = 2261.926
$ws.Range("I89").Value = 2358.7727
$ws.Range("K89").Value = 11793.8635
$ws.Range("M89").Value = -6177.863499999999

$ws.Range("H94").Value = 461.25
$ws.Range("I94").Value = 348.5
$ws.Range("J94").Value = 574
$ws.Range("K94").Value = 348.5
$ws.Range("L94").Value = 574
$ws.Range("N94").Value = -1476
$ws.Range("M94").Value = 102.5

$ws.Range("H107").Value = 620.25
$ws.Range("I107").Value = 620.25
$ws.Range("K107").Value = 620.25
$ws.Range("M107").Value = 1299.75

$ws.Range("H134").Value = 2947
$ws.Range("I134").Value = 3067.6667
$ws.Range("J134").Value = 1499
$ws.Range("K134").Value = 9203.000100000001
$ws.Range("L134").Value = 4497
$ws.Range("M134").Value = -6668.000100000001
$ws.Range("N134").Value = -9567

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1499.5

$ws.Range("H22").Value = 53500.5
$ws.Range("J22").Value = 53500.5
$ws.Range("L22").Value = 53500.5
$ws.Range("N22").Value = -54200.5

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null

$ws.Range("H99").Value = 2341.4285
$ws.Range("I99").Value = 2248.3333
$ws.Range("K99").Value = 2248.3333
$ws.Range("M99").Value = -750.3332999999998

$ws.Range("H126").Value = 2341.4285
$ws.Range("I126").Value = 2248.3333
$ws.Range("K126").Value = 6744.999899999999
$ws.Range("M126").Value = -4274.999899999999

$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1595.75
$ws.Range("I11").Value = 189.5
$ws.Range("K11").Value = 568.5
$ws.Range("M11").Value = -428.5

$ws.Range("H26").Value = 1546.1515
$ws.Range("J26").Value = 1655.8572
$ws.Range("L26").Value = 4967.571599999999
$ws.Range("N26").Value = -5543.571599999999

$ws.Range("H41").Value = 2197.25
$ws.Range("J41").Value = 4244.5
$ws.Range("L41").Value = 12733.5
$ws.Range("N41").Value = -13409.5

$ws.Range("H44").Value = 300
$ws.Range("I44").Value = 300
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 900
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -502
$ws.Range("N44").Value = $null

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2241.2856
$ws.Range("I80").Value = 2525
$ws.Range("J80").Value = 1863
$ws.Range("K80").Value = 2525
$ws.Range("L80").Value = 1863
$ws.Range("M80").Value = -1527
$ws.Range("N80").Value = -3859

$ws.Range("H83").Value = 2241.2856
$ws.Range("I83").Value = 2525
$ws.Range("J83").Value = 1863
$ws.Range("K83").Value = 12625
$ws.Range("L83").Value = 9315
$ws.Range("M83").Value = -7633
$ws.Range("N83").Value = -19299

$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -46134

$ws.Range("H113").Value = 2006.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2006.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2006.5
$ws.Range("N113").Value = -6346.5
$ws.Range("M113").Value = $null

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 14999.5
$ws.Range("J2").Value = 14999.5
$ws.Range("L2").Value = 14999.5
$ws.Range("N2").Value = -15223.5

$ws.Range("H21").Value = 1000
$ws.Range("J21").Value = 1000
$ws.Range("L21").Value = 1000
$ws.Range("N21").Value = -1348

$ws.Range("H55").Value = 869.3077
$ws.Range("I55").Value = 700.2857
$ws.Range("K55").Value = 700.2857
$ws.Range("M55").Value = -527.2857

$ws.Range("H132").Value = 3659.5715
$ws.Range("I132").Value = 3649.4546
$ws.Range("J132").Value = 3696.6667
$ws.Range("K132").Value = 10948.3638
$ws.Range("L132").Value = 11090.0001
$ws.Range("M132").Value = -8418.363799999999
$ws.Range("N132").Value = -16150.0001

$ws.Range("H136").Value = 1450
$ws.Range("I136").Value = 1450
$ws.Range("K136").Value = 4350
$ws.Range("M136").Value = -1800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 50022496
$ws.Range("I2").Value = 50022496
$ws.Range("K2").Value = 50022496
$ws.Range("M2").Value = -50022384

$ws.Range("H4").Value = 3999.5
$ws.Range("J4").Value = 3999.5
$ws.Range("L4").Value = 3999.5
$ws.Range("N4").Value = -4225.5

$ws.Range("H122").Value = 3061.75
$ws.Range("I122").Value = 2984.8572
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 8954.571599999999
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -6504.571599999999
$ws.Range("N122").Value = -15700

$ws.Range("H133").Value = 53333.332
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120

$ws.Range("H136").Value = 3686.524
$ws.Range("I136").Value = 3819.6086
$ws.Range("K136").Value = 11458.8258
$ws.Range("M136").Value = -8908.825800000001
